$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top; existing data (rows 1-47) shifts down to rows 3-49
$ws.Range("A1:C2").Insert()

# New row 1: generic "Unnamed: N" headers (bold, centered, top-aligned, thin box border)
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "Unnamed: 1"
$ws.Range("C1").Value = "Unnamed: 2"

# New row 2: real column headers
$ws.Range("A2").Value = "municipio"
$ws.Range("B2").Value = "Casos"
$ws.Range("C2").Value = "Óbitos"

# Style row 1 like a header band: bold font, centered horizontally, top vertical
# alignment, thin border all around each cell
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
